# Refresh the cached "datetimeFigureOut" date placeholder text across the
# Slide Master and every Slide Layout (Insert > Header & Footer style
# update), e.g. 23/11/2020 -> 02/06/2023.

$p = $ppt.ActivePresentation
$newDate = "02/06/2023"
$ppPlaceholderDate = 16

function Update-IfDatePlaceholder($shape) {
    if ($shape.HasTextFrame -eq -1) {
        $isDate = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if (-not $isDate -and $shape.Name -like "Date Placeholder*") {
            $isDate = $true
        }
        if ($isDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    Update-IfDatePlaceholder($master.Shapes.Item($j))
}

# Every Slide Layout off the master
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        Update-IfDatePlaceholder($layout.Shapes.Item($j))
    }
}
